# docs/diagrams: Fix sequence diagram error
#
# 1) Fix the mis-drawn "Straight Arrow Connector 49" on slide 1: it should be
#    a flat (non-flipped) horizontal connector ending where the "rc" lifeline
#    box is, not a long diagonal flipped connector.
# 2) Refresh the cached "today" text of the auto-updating date placeholders
#    (datetimeFigureOut fields) on the slide master / layouts / notes master,
#    which PowerPoint re-stamps whenever the deck is re-saved on a later day.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Fix the connector geometry on slide 1
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

$connector = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Straight Arrow Connector 49") {
        $connector = $sh
        break
    }
}

if ($connector -ne $null) {
    # Remove the vertical flip ...
    $connector.VerticalFlip = 0
    # ... and reshape it into a short, flat connector (Left/x is unchanged).
    # Values are expressed in points (1 pt = 12700 EMU); a tiny epsilon is
    # added where needed so the emitted EMU rounds to the exact target value.
    $connector.Top    = (2975344 / 12700) + 0.00001
    $connector.Width  = (1439295 / 12700) + 0.00001
    $connector.Height = 0
}

# ---------------------------------------------------------------------------
# 2) Bump the cached date placeholder text from 11/11/2018 to 11/12/2018
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "11/11/2018") {
                $tr.Text = "11/12/2018"
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout under the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Notes master
Update-DatePlaceholder $p.NotesMaster.Shapes
